$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcRow = 8
$dstRow = 9

# Copy formatting (including number formats/styles) from the row above
# so the new row matches existing style conventions (e.g. date column G).
$ws.Range($ws.Cells.Item($srcRow, 1), $ws.Cells.Item($srcRow, 8)).Copy($ws.Cells.Item($dstRow, 1))

$ws.Cells.Item($dstRow, 1).Value = 9758.1299999999992
$ws.Cells.Item($dstRow, 2).Value = 9815.06
$ws.Cells.Item($dstRow, 3).Value = 307.20999999999998
$ws.Cells.Item($dstRow, 4).Value = 305.42
$ws.Cells.Item($dstRow, 5).Value = $false
$ws.Cells.Item($dstRow, 6).Value = -0.57999999999999996
$ws.Cells.Item($dstRow, 7).Value = 42609.488958333335
$ws.Cells.Item($dstRow, 8).Value = $false
